$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-10 11:51:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-10 11:51:16"
